# Update crypto price/volume figures to the latest scraped values.
# Values are stored as literal text in the sheet (e.g. "278.29", "0.90%"),
# so each cell is forced to text format before assignment, then the style is
# reset back to Normal so no stray formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($Sheet, $Address, $Text) {
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "278.29"
Set-TextValue $ws "E2" "0.90%"
Set-TextValue $ws "E3" "2.03%"
Set-TextValue $ws "D4" "4.873"
Set-TextValue $ws "E4" "0.02%"
Set-TextValue $ws "D5" "0.06430"
Set-TextValue $ws "E5" "1.60%"
Set-TextValue $ws "D6" "7.003"
Set-TextValue $ws "E6" "1.38%"
Set-TextValue $ws "D7" "1.205"
Set-TextValue $ws "E7" "-6.59%"
Set-TextValue $ws "D8" "0.8893"
Set-TextValue $ws "E8" "1.62%"
Set-TextValue $ws "D9" "0.1552"
Set-TextValue $ws "E9" "-1.89%"
Set-TextValue $ws "E10" "1.73%"
Set-TextValue $ws "D11" "0.07507"
Set-TextValue $ws "E11" "0.51%"
Set-TextValue $ws "D12" "0.02887"
Set-TextValue $ws "E12" "-2.30%"
Set-TextValue $ws "D13" "0.08968"
Set-TextValue $ws "E13" "-1.03%"
Set-TextValue $ws "D14" "0.001570"
Set-TextValue $ws "E14" "-0.86%"
Set-TextValue $ws "D15" "0.0006375"
Set-TextValue $ws "E15" "0.95%"
Set-TextValue $ws "D16" "0.006111"
Set-TextValue $ws "E16" "1.62%"
Set-TextValue $ws "D17" "3.479"
Set-TextValue $ws "E17" "1.03%"
Set-TextValue $ws "D18" "3.306"
Set-TextValue $ws "E18" "-0.32%"
Set-TextValue $ws "E19" "-2.87%"
Set-TextValue $ws "E20" "1.11%"
Set-TextValue $ws "E21" "1.33%"
Set-TextValue $ws "D22" "3.916"
Set-TextValue $ws "E22" "0.21%"
Set-TextValue $ws "D23" "0.04403"
Set-TextValue $ws "E23" "0.72%"
Set-TextValue $ws "E24" "8.71%"
Set-TextValue $ws "D25" "0.001175"
Set-TextValue $ws "E25" "0.31%"
Set-TextValue $ws "D26" "0.003878"
Set-TextValue $ws "E26" "-7.94%"
Set-TextValue $ws "E28" "-1.70%"
Set-TextValue $ws "E29" "1.74%"
Set-TextValue $ws "D40" "0.04125"
Set-TextValue $ws "E40" "0.77%"
Set-TextValue $ws "D41" "0.006802"
Set-TextValue $ws "E41" "-2.59%"
Set-TextValue $ws "E42" "0.10%"
Set-TextValue $ws "D43" "0.001920"
Set-TextValue $ws "E43" "-9.89%"
Set-TextValue $ws "D44" "0.01185"
Set-TextValue $ws "E44" "9.49%"
Set-TextValue $ws "D45" "0.00005328"
Set-TextValue $ws "E45" "0.45%"
Set-TextValue $ws "E46" "13.03%"
Set-TextValue $ws "D47" "0.01852"
Set-TextValue $ws "E47" "-7.40%"
